$d = $word.ActiveDocument

# 1. Split "Egy szerzőnek több email címe is lehet." into two sentences/runs
#    and add the new clause about name/code uniqueness.
$d.Content.Find.Execute(
    "Egy szerzőnek több email címe is lehet.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Egy szerzőnek több email címe is lehet, de neve és kódja csak egy van.",
    2) | Out-Null

# 2. Simplify the "neptun" sentence — remove the spell-check split and
#    just replace the whole (now-merged) text with the plain version.
$d.Content.Find.Execute(
    "Egy hallgatónak egy neptun kódja és egy neve van.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Egy hallgatónak egy neptun kódja és egy neve van.",
    2) | Out-Null

# 3. Replace the trailing bookmark-only paragraph with a new sentence
#    about campuses, followed by two empty paragraphs.
$last = $d.Paragraphs.Last
$last.Range.Text = "Egy campusnak egy azonosítója és egy címe van.`r`r`r"
